# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 76 of Sheet1, pushing the
# existing rows 76-152 down to 77-153 (dimension grows from A1:T152 to
# A1:T153).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 76; Excel shifts everything below it down
# by one and the sheet dimension is extended automatically.
$ws.Rows("76:76").Insert()

# Populate the new row with the latest observation (constant columns
# copied from the rest of the block; Fecha/Calidad/Volumen/Precios/
# Origen/Precio per-Kg are the new data point).
$ws.Range("A76").Value = 9
$ws.Range("B76").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C76").Value = "Metropolitana"
$ws.Range("D76").Value = 44994
$ws.Range("E76").Value = 13
$ws.Range("F76").Value = "Fruta"
$ws.Range("G76").Value = 100101
$ws.Range("H76").Value = "Berries"
$ws.Range("I76").Value = 100101004
$ws.Range("J76").Value = "Frambuesa"
$ws.Range("K76").Value = "Sin especificar"
$ws.Range("L76").Value = "Primera"
$ws.Range("M76").Value = 580
$ws.Range("N76").Value = 6000
$ws.Range("O76").Value = 6500
$ws.Range("P76").Value = 6241
$ws.Range("Q76").Value = "$/bandeja 2 kilos"
$ws.Range("R76").Value = "Provincia de Curicó"
$ws.Range("S76").Value = 3120
$ws.Range("T76").Value = 2
